$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 165 (shifts "Rohan Mudras" and everything below it down by one row)
$ws.Rows.Item(165).Insert()

# Populate the new row 165 with the "Riya Arora" entry
$ws.Range("A165").Value = "Riya Arora"

# Force the phone number into B165 as text (not a number) to match column formatting
$ws.Range("B165").NumberFormat = "@"
$ws.Range("B165").Value = "9854615863"
$ws.Range("B165").ClearFormats()

$ws.Range("C165").Formula = '=HYPERLINK("https://qrcode-2-production.up.railway.app/qr/Riya_Arora_1873ebb2.png","https://qrcode-2-production.up.railway.app/qr/Riya_Arora_1873ebb2.png")'
